# Update "想去人数" (column F) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value for column F
$updates = @{
    2  = 2028
    4  = 125
    5  = 42
    7  = 1687
    9  = 680
    10 = 367
    19 = 3907
    22 = 441
    23 = 367
    24 = 828
    25 = 555
    26 = 358
    28 = 1702
    29 = 17
    31 = 169
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
